$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-13, columns A-T
# Values written column-by-column (A first, then B, C, ... ) to match
# the source data layout / shared-string insertion order.

# Column A
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "ECs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "FAPs"
$ws.Range("A9").Value = "FAPs"
$ws.Range("A10").Value = "sCs"
$ws.Range("A11").Value = "sCs"
$ws.Range("A12").Value = "sCs"
$ws.Range("A13").Value = "sCs"

# Column B
$ws.Range("B2").Value = "Col4a5"
$ws.Range("B3").Value = "Col4a5"
$ws.Range("B4").Value = "Col4a5"
$ws.Range("B5").Value = "Col4a5"
$ws.Range("B6").Value = "Col4a5"
$ws.Range("B7").Value = "Col4a5"
$ws.Range("B8").Value = "Col4a5"
$ws.Range("B9").Value = "Col4a5"
$ws.Range("B10").Value = "Col4a5"
$ws.Range("B11").Value = "Col4a5"
$ws.Range("B12").Value = "Col4a5"
$ws.Range("B13").Value = "Col4a5"

# Column C
$ws.Range("C2").Value = "Cd93"
$ws.Range("C3").Value = "Cd93"
$ws.Range("C4").Value = "Cd93"
$ws.Range("C5").Value = "Cd93"
$ws.Range("C6").Value = "Cd93"
$ws.Range("C7").Value = "Cd93"
$ws.Range("C8").Value = "Cd93"
$ws.Range("C9").Value = "Cd93"
$ws.Range("C10").Value = "Cd93"
$ws.Range("C11").Value = "Cd93"
$ws.Range("C12").Value = "Cd93"
$ws.Range("C13").Value = "Cd93"

# Column D
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "M2"
$ws.Range("D5").Value = "sCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("D8").Value = "M2"
$ws.Range("D9").Value = "sCs"
$ws.Range("D10").Value = "ECs"
$ws.Range("D11").Value = "FAPs"
$ws.Range("D12").Value = "M2"
$ws.Range("D13").Value = "sCs"

# Column E
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("E10").Value = 3
$ws.Range("E11").Value = 3
$ws.Range("E12").Value = 3
$ws.Range("E13").Value = 3

# Column F
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1

# Column G
$ws.Range("G2").Value = 0.037247
$ws.Range("G3").Value = 0.037247
$ws.Range("G4").Value = 0.037247
$ws.Range("G5").Value = 0.037247
$ws.Range("G6").Value = 4.613664666666667
$ws.Range("G7").Value = 4.613664666666667
$ws.Range("G8").Value = 4.613664666666667
$ws.Range("G9").Value = 4.613664666666667
$ws.Range("G10").Value = 2.364937
$ws.Range("G11").Value = 2.364937
$ws.Range("G12").Value = 2.364937
$ws.Range("G13").Value = 2.364937

# Column H
$ws.Range("H2").Value = 0.111741
$ws.Range("H3").Value = 0.111741
$ws.Range("H4").Value = 0.111741
$ws.Range("H5").Value = 0.111741
$ws.Range("H6").Value = 13.840994
$ws.Range("H7").Value = 13.840994
$ws.Range("H8").Value = 13.840994
$ws.Range("H9").Value = 13.840994
$ws.Range("H10").Value = 7.094811
$ws.Range("H11").Value = 7.094811
$ws.Range("H12").Value = 7.094811
$ws.Range("H13").Value = 7.094811

# Column I
$ws.Range("I2").Value = 0.005308979963744942
$ws.Range("I3").Value = 0.005308979963744942
$ws.Range("I4").Value = 0.005308979963744942
$ws.Range("I5").Value = 0.005308979963744942
$ws.Range("I6").Value = 0.6576060696102054
$ws.Range("I7").Value = 0.6576060696102054
$ws.Range("I8").Value = 0.6576060696102054
$ws.Range("I9").Value = 0.6576060696102054
$ws.Range("I10").Value = 0.3370849504260496
$ws.Range("I11").Value = 0.3370849504260496
$ws.Range("I12").Value = 0.3370849504260496
$ws.Range("I13").Value = 0.3370849504260496

# Column J
$ws.Range("J2").Value = 0.005308979963744942
$ws.Range("J3").Value = 0.005308979963744942
$ws.Range("J4").Value = 0.005308979963744942
$ws.Range("J5").Value = 0.005308979963744942
$ws.Range("J6").Value = 0.6576060696102054
$ws.Range("J7").Value = 0.6576060696102054
$ws.Range("J8").Value = 0.6576060696102054
$ws.Range("J9").Value = 0.6576060696102054
$ws.Range("J10").Value = 0.3370849504260496
$ws.Range("J11").Value = 0.3370849504260496
$ws.Range("J12").Value = 0.3370849504260496
$ws.Range("J13").Value = 0.3370849504260496

# Column K
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 2
$ws.Range("K8").Value = 3
$ws.Range("K9").Value = 3
$ws.Range("K10").Value = 3
$ws.Range("K11").Value = 2
$ws.Range("K12").Value = 3
$ws.Range("K13").Value = 3

# Column L
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("L12").Value = 1
$ws.Range("L13").Value = 1

# Column M
$ws.Range("M2").Value = 135.955556
$ws.Range("M3").Value = 0.449122
$ws.Range("M4").Value = 56.38366533333333
$ws.Range("M5").Value = 2.896484
$ws.Range("M6").Value = 135.955556
$ws.Range("M7").Value = 0.449122
$ws.Range("M8").Value = 56.38366533333333
$ws.Range("M9").Value = 2.896484
$ws.Range("M10").Value = 135.955556
$ws.Range("M11").Value = 0.449122
$ws.Range("M12").Value = 56.38366533333333
$ws.Range("M13").Value = 2.896484

# Column N
$ws.Range("N2").Value = 407.866668
$ws.Range("N3").Value = 1.347366
$ws.Range("N4").Value = 169.150996
$ws.Range("N5").Value = 8.689452
$ws.Range("N6").Value = 407.866668
$ws.Range("N7").Value = 1.347366
$ws.Range("N8").Value = 169.150996
$ws.Range("N9").Value = 8.689452
$ws.Range("N10").Value = 407.866668
$ws.Range("N11").Value = 1.347366
$ws.Range("N12").Value = 169.150996
$ws.Range("N13").Value = 8.689452

# Column O
$ws.Range("O2").Value = 0.6947679994035034
$ws.Range("O3").Value = 0.002295129398228494
$ws.Range("O4").Value = 0.2881350899898248
$ws.Range("O5").Value = 0.01480178120844327
$ws.Range("O6").Value = 0.6947679994035034
$ws.Range("O7").Value = 0.002295129398228494
$ws.Range("O8").Value = 0.2881350899898248
$ws.Range("O9").Value = 0.01480178120844327
$ws.Range("O10").Value = 0.6947679994035034
$ws.Range("O11").Value = 0.002295129398228494
$ws.Range("O12").Value = 0.2881350899898248
$ws.Range("O13").Value = 0.01480178120844327

# Column P
$ws.Range("P2").Value = 0.6947679994035034
$ws.Range("P3").Value = 0.002295129398228494
$ws.Range("P4").Value = 0.2881350899898248
$ws.Range("P5").Value = 0.01480178120844327
$ws.Range("P6").Value = 0.6947679994035034
$ws.Range("P7").Value = 0.002295129398228494
$ws.Range("P8").Value = 0.2881350899898248
$ws.Range("P9").Value = 0.01480178120844327
$ws.Range("P10").Value = 0.6947679994035034
$ws.Range("P11").Value = 0.002295129398228494
$ws.Range("P12").Value = 0.2881350899898248
$ws.Range("P13").Value = 0.01480178120844327

# Column Q
$ws.Range("Q2").Value = 5.063936594332
$ws.Range("Q3").Value = 0.016728447134
$ws.Range("Q4").Value = 2.100122382670667
$ws.Range("Q5").Value = 0.107885339548
$ws.Range("Q6").Value = 627.2533449542213
$ws.Range("Q7").Value = 2.072098302422667
$ws.Range("Q8").Value = 260.1353245255582
$ws.Range("Q9").Value = 13.36340588836533
$ws.Range("Q10").Value = 321.526324739972
$ws.Range("Q11").Value = 1.062145235314
$ws.Range("Q12").Value = 133.3438163424173
$ws.Range("Q13").Value = 6.850002181507999

# Column R
$ws.Range("R2").Value = 45.575429348988
$ws.Range("R3").Value = 0.150556024206
$ws.Range("R4").Value = 18.901101444036
$ws.Range("R5").Value = 0.9709680559319999
$ws.Range("R6").Value = 5645.280104587992
$ws.Range("R7").Value = 18.648884721804
$ws.Range("R8").Value = 2341.217920730024
$ws.Range("R9").Value = 120.270652995288
$ws.Range("R10").Value = 2893.736922659748
$ws.Range("R11").Value = 9.559307117826
$ws.Range("R12").Value = 1200.094347081756
$ws.Range("R13").Value = 61.65001963357199

# Column S
$ws.Range("S2").Value = 0.003688509388284357
$ws.Range("S3").Value = [double]"1.218479598939706E-05"
$ws.Range("S4").Value = 0.001529703419607826
$ws.Range("S5").Value = [double]"7.858235986336173E-05"
$ws.Range("S6").Value = 0.4568836533786834
$ws.Range("S7").Value = 0.001509291022815876
$ws.Range("S8").Value = 0.1894793840449916
$ws.Range("S9").Value = 0.009733741163714577
$ws.Range("S10").Value = 0.2341958366365356
$ws.Range("S11").Value = 0.000773653579423221
$ws.Range("S12").Value = 0.09712600252522546
$ws.Range("S13").Value = 0.004989457684865334

# Column T
$ws.Range("T2").Value = 0.003688509388284357
$ws.Range("T3").Value = [double]"1.218479598939706E-05"
$ws.Range("T4").Value = 0.001529703419607826
$ws.Range("T5").Value = [double]"7.858235986336173E-05"
$ws.Range("T6").Value = 0.4568836533786834
$ws.Range("T7").Value = 0.001509291022815876
$ws.Range("T8").Value = 0.1894793840449915
$ws.Range("T9").Value = 0.009733741163714577
$ws.Range("T10").Value = 0.2341958366365356
$ws.Range("T11").Value = 0.0007736535794232209
$ws.Range("T12").Value = 0.09712600252522544
$ws.Range("T13").Value = 0.004989457684865334
